$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "5) RMSE" sheet - update the "complete cases" RMSE values
# ------------------------------------------------------------------
$wsRmse = $wb.Worksheets.Item("5) RMSE")

$wsRmse.Range("C3").Value = 1.52

$wsRmse.Range("B4").Value = 0.68
$wsRmse.Range("C4").Value = 1.75

$wsRmse.Range("B5").Value = 2.39
$wsRmse.Range("C5").Value = 1.18

$wsRmse.Range("B6").Value = 0.47
$wsRmse.Range("C6").Value = 0.76

$wsRmse.Range("B7").Value = 0.3
$wsRmse.Range("C7").Value = 0.99

# ------------------------------------------------------------------
# 2) "4) Mass balances" sheet - re-enter the shared formula in E13:F13
#    (keeps the same formula text / computed values; tidies the
#    shared-formula bookkeeping for that row)
# ------------------------------------------------------------------
$wsMass = $wb.Worksheets.Item("4) Mass balances")
$wsMass.Range("E13:F13").Formula = '=E4/$G$4'

# ------------------------------------------------------------------
# 3) Active sheet / tab selection moves from
#    "3) Built model with equations" to "5) RMSE", with the RMSE
#    sheet's selection landing on E7.
# ------------------------------------------------------------------
$wsRmse.Activate()
$wsRmse.Range("E7").Select()
